$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-15 12:07:41", 0.0004),
    @("2023-12-15 12:07:56", 0.0008),
    @("2023-12-15 12:08:29", 0.0034),
    @("2023-12-15 12:08:34", 0.0002),
    @("2023-12-15 12:08:44", 0.0004),
    @("2023-12-15 12:08:53", 0.0002)
)

$startRow = 355
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
